$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The column D header was renamed from "coral_id" to "specie_id".
$ws.Range("D1").Value = "specie_id"

# Reflect the active cell / selection recorded in the sheet view.
$ws.Range("D1").Select()
